$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, $Text)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "29.858.83"
$ws.Range("E2").Value = "  +0.43%  "

# Row 3
Set-TextValue $ws.Range("D3") "1.893.43"
$ws.Range("E3").Value = "  +0.36%  "

# Row 4
$ws.Range("E4").Value = "  +0.06%  "

# Row 5
Set-TextValue $ws.Range("D5") "0.7807"
$ws.Range("E5").Value = "  -1.31%  "

# Row 6
Set-TextValue $ws.Range("D6") "243.64"
$ws.Range("E6").Value = "  +0.86%  "

# Row 7
$ws.Range("E7").Value = "  +0.00%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.3134"
$ws.Range("E8").Value = "  -0.91%  "

# Row 9
Set-TextValue $ws.Range("D9") "25.70"
$ws.Range("E9").Value = "  +1.54%  "

# Row 10
Set-TextValue $ws.Range("D10") "0.07273"
$ws.Range("E10").Value = "  +4.00%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.08092"
$ws.Range("E11").Value = "  +0.59%  "

# Row 12
Set-TextValue $ws.Range("D12") "0.7726"
$ws.Range("E12").Value = "  +1.01%  "

# Row 13
Set-TextValue $ws.Range("D13") "5.463"
$ws.Range("E13").Value = "  +3.33%  "

# Row 14
Set-TextValue $ws.Range("D14") "1.884.70"
$ws.Range("E14").Value = "  +0.07%  "

# Row 15
Set-TextValue $ws.Range("D15") "93.88"
$ws.Range("E15").Value = "  +2.18%  "

# Row 16
$ws.Range("E16").Value = "  +4.99%  "

# Row 17
Set-TextValue $ws.Range("D17") "29.868.99"
$ws.Range("E17").Value = "  +0.42%  "

# Row 18
Set-TextValue $ws.Range("D18") "13.93"
$ws.Range("E18").Value = "  +0.67%  "

# Row 19
Set-TextValue $ws.Range("D19") "246.58"
$ws.Range("E19").Value = "  +1.42%  "

# Row 20
Set-TextValue $ws.Range("D20") "0.000007806"
$ws.Range("E20").Value = "  +1.58%  "

# Row 21
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Range("D21") "1.000"
$ws.Range("E21").Value = "  +0.01%  "

# Row 22
$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue $ws.Range("D22") "2.147.32"
$ws.Range("E22").Value = "  +0.74%  "

# Row 23
Set-TextValue $ws.Range("D23") "8.114"
$ws.Range("E23").Value = "  -0.54%  "

# Row 24
Set-TextValue $ws.Range("D24") "1.001"
$ws.Range("E24").Value = "  +0.02%  "

# Row 25
$ws.Range("E25").Value = "  -4.45%  "

# Row 26
Set-TextValue $ws.Range("D26") "9.436"
$ws.Range("E26").Value = "  +1.71%  "

# Row 27
Set-TextValue $ws.Range("D27") "164.33"
$ws.Range("E27").Value = "  -0.45%  "

# Row 28
Set-TextValue $ws.Range("D28") "18.73"
$ws.Range("E28").Value = "  +0.62%  "

# Row 29
Set-TextValue $ws.Range("D29") "2.019"
$ws.Range("E29").Value = "  -1.15%  "

# Row 30
Set-TextValue $ws.Range("D30") "1.440"
$ws.Range("E30").Value = "  +2.90%  "

# Row 31
$ws.Range("E31").Value = "  +0.75%  "

# Row 32
Set-TextValue $ws.Range("D32") "4.477"
$ws.Range("E32").Value = "  +2.10%  "

# Row 33
$ws.Range("E33").Value = "  -2.10%  "

# Row 34
Set-TextValue $ws.Range("D34") "4.063"
$ws.Range("E34").Value = "  +0.57%  "

# Row 35
Set-TextValue $ws.Range("D35") "1.241"
$ws.Range("E35").Value = "  -1.42%  "

# Row 36
Set-TextValue $ws.Range("D36") "0.7526"
$ws.Range("E36").Value = "  +2.56%  "

# Row 37
$ws.Range("E37").Value = "  +0.62%  "

# Row 38
Set-TextValue $ws.Range("D38") "2.681"
$ws.Range("E38").Value = "  +1.78%  "

# Row 39
Set-TextValue $ws.Range("D39") "0.01935"
$ws.Range("E39").Value = "  +1.48%  "

# Row 40
Set-TextValue $ws.Range("D40") "2.801"
$ws.Range("E40").Value = "  +1.23%  "

# Row 41
Set-TextValue $ws.Range("D41") "1.137.31"
$ws.Range("E41").Value = "  +11.38%  "

# Row 42
Set-TextValue $ws.Range("D42") "0.4460"
$ws.Range("E42").Value = "  +1.63%  "

# Row 43
Set-TextValue $ws.Range("D43") "74.16"
$ws.Range("E43").Value = "  +2.71%  "

# Row 44
Set-TextValue $ws.Range("D44") "5.964"
$ws.Range("E44").Value = "  +2.45%  "

# Row 45
Set-TextValue $ws.Range("D45") "0.8527"
$ws.Range("E45").Value = "  +2.10%  "

# Row 46
Set-TextValue $ws.Range("D46") "1.000"
$ws.Range("E46").Value = "  -0.01%  "

# Row 47
Set-TextValue $ws.Range("D47") "1.886"
$ws.Range("E47").Value = "  +1.83%  "

# Row 48
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws.Range("D48") "102.45"
$ws.Range("E48").Value = "  +0.21%  "

# Row 49
$ws.Range("B49").Value = "SynthetixNetwork"
$ws.Range("C49").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
Set-TextValue $ws.Range("D49") "3.113"
$ws.Range("E49").Value = "  +7.68%  "

# Row 50
Set-TextValue $ws.Range("D50") "7.530"
$ws.Range("E50").Value = "  +1.81%  "

# Row 51
Set-TextValue $ws.Range("D51") "9.738"
$ws.Range("E51").Value = "  -1.40%  "
